# The presentation originally carries the "Integral" design (Red Violet
# colour scheme) on ppt/theme/theme1.xml, which is the theme bound to the
# slide master / all slide layouts / the whole deck. The commit swaps the
# active colour scheme over to the standard "Office" palette (the colours
# that used to live, unused, in ppt/theme/theme2.xml - the Notes Master's
# theme).
#
# PowerPoint's object model exposes the twelve theme colour slots (dk1,
# lt1, dk2, lt2, accent1-6, hlink, folHlink) via Slide.ThemeColorScheme,
# in that exact order, each as a settable .RGB value. Updating them here
# rewrites <a:clrScheme> inside the theme part that backs the slide
# master (ppt/theme/theme1.xml), which is precisely the visual change the
# commit makes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target values are the standard Office theme colours, expressed as
# VBA-style RGB() decimals (R + G*256 + B*65536), in ThemeColorScheme
# slot order.
$officeColors = @(
    0        # 1  dk1      000000
    16777215 # 2  lt1      FFFFFF
    6968388  # 3  dk2      44546A
    15132391 # 4  lt2      E7E6E6
    13998939 # 5  accent1  5B9BD5
    3243501  # 6  accent2  ED7D31
    10855845 # 7  accent3  A5A5A5
    49407    # 8  accent4  FFC000
    12874308 # 9  accent5  4472C4
    4697456  # 10 accent6  70AD47
    12673797 # 11 hlink    0563C1
    7491477  # 12 folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
